# (OPR)Import Payroll Bonus.xlsx - update payroll and system
# Applies:
#  - sheet2 ("description"): replace the old "item number" note with the
#    longer "numbers only" note, add two new help columns H:I (merged header
#    + two descriptive cells), resize H:I, select B9
#  - sheet1: select B13 (becomes the active sheet/tab)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- sheet2: rebuild the shared-string table in the right order -----------
# 1) Drop the old note text so its shared-string slot is freed/compacted.
$ws2.Range("B1").Value = ""

# 2) Add the two brand-new help strings (H2/I2) BEFORE the replacement note,
#    so they land at the matching shared-string indices once saved.
$ws2.Range("H2").Value = "เป็นค่าว่างได้"
$ws2.Range("I2").Value = "เพิ่มข้อมูล ต้องไม่เป็นค่าว่าง"

# 3) Now write the new (longer) note text back into B1.
$ws2.Range("B1").Value = "เลขที่  (null)  ตัวเลขเท่านั้น"

# --- sheet2: styles for the new H2/I2 cells (reuse existing named styles) --
$ws2.Range("B2").Copy()
$ws2.Range("H2").PasteSpecial(-4122)
$ws2.Range("C2").Copy()
$ws2.Range("I2").PasteSpecial(-4122)

# --- sheet2: new merged header cell H1:I1 (centered, no fill/border) -------
$ws2.Range("H1:I1").VerticalAlignment = -4107
$ws2.Range("H1:I1").HorizontalAlignment = -4108
$ws2.Range("H1:I1").Merge()

# --- sheet2: widen the two new columns -------------------------------------
$ws2.Columns("H:I").ColumnWidth = 24.8

# --- sheet2: selection moves to B9 (no longer the active tab) -------------
$ws2.Range("B9").Select()

# --- sheet1: becomes the active sheet/tab, selection moves to B13 ---------
$ws1.Range("B13").Select()
